$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 contains the user details record; update the values that changed.
$ws.Range("A2").Value = "WpfUg383"      # Client Id
$ws.Range("B2").Value = 23080282        # Candidate ID
$ws.Range("C2").Value = "cijgcsa64"     # User Name
$ws.Range("D2").Value = "b8&HfG3%"      # Exam Password
$ws.Range("F2").Value = "OtwOukde"      # First Name
$ws.Range("G2").Value = "JVWw"          # Last Name
